$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to restore style after
# forcing text NumberFormat so numeric-looking strings are not auto-converted
# to Excel numbers (preserving the original inlineStr "text" cell semantics).
$refStyle = $ws.Range("F2").Style

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $refStyle
}

Set-TextValue "D2" "258.48"
Set-TextValue "E2" "0.73%"
Set-TextValue "E3" "-1.30%"
Set-TextValue "D4" "4.649"
Set-TextValue "E4" "1.97%"
Set-TextValue "D5" "0.05982"
Set-TextValue "E5" "1.50%"
Set-TextValue "D6" "6.645"
Set-TextValue "D7" "0.8563"
Set-TextValue "E7" "-0.27%"
Set-TextValue "D8" "0.9201"
Set-TextValue "E8" "-0.96%"
Set-TextValue "D9" "0.1387"
Set-TextValue "E9" "-1.53%"
Set-TextValue "D10" "0.04896"
Set-TextValue "E10" "34.99%"
Set-TextValue "D11" "0.07024"
Set-TextValue "E11" "-1.06%"
Set-TextValue "D12" "0.03047"
Set-TextValue "E12" "-5.77%"
Set-TextValue "D13" "0.09122"
Set-TextValue "E13" "-0.96%"
Set-TextValue "D14" "0.001525"
Set-TextValue "E14" "-1.19%"
Set-TextValue "D15" "0.0006069"
Set-TextValue "E15" "0.13%"
Set-TextValue "D16" "0.006179"
Set-TextValue "E16" "1.44%"
Set-TextValue "D17" "3.447"
Set-TextValue "E17" "-1.95%"
Set-TextValue "D18" "3.147"
Set-TextValue "E18" "-1.35%"
Set-TextValue "D19" "2.191"
Set-TextValue "E19" "-0.50%"
Set-TextValue "D20" "0.3110"
Set-TextValue "E20" "1.76%"
Set-TextValue "E21" "0.89%"
Set-TextValue "D22" "4.047"
Set-TextValue "E22" "5.26%"
Set-TextValue "D23" "0.04214"
Set-TextValue "E23" "0.07%"
Set-TextValue "E24" "-0.60%"
Set-TextValue "D25" "0.004023"
Set-TextValue "E25" "-5.96%"
Set-TextValue "E26" "-0.05%"
Set-TextValue "E27" "13.30%"
Set-TextValue "D40" "0.03825"
Set-TextValue "E40" "-0.19%"
Set-TextValue "D41" "0.1111"
Set-TextValue "E41" "1.13%"
Set-TextValue "D42" "0.003768"
Set-TextValue "E42" "-39.43%"
Set-TextValue "D43" "0.002429"
Set-TextValue "E43" "10.39%"
Set-TextValue "D44" "0.01520"
Set-TextValue "E44" "33.71%"
Set-TextValue "D45" "0.00005099"
Set-TextValue "E45" "-6.10%"
Set-TextValue "E46" "-0.05%"
Set-TextValue "E47" "-43.24%"
Set-TextValue "D48" "0.1072"
Set-TextValue "E48" "2.30%"
Set-TextValue "E49" "-0.05%"
Set-TextValue "E50" "-0.05%"
